$d = $word.ActiveDocument

function Replace-ExactRange([string]$searchText, [string]$replaceText) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $r.Text = $replaceText
    return $r
}

# 1) GET 'normal_posts' -> GET 'posts'
Replace-ExactRange "GET ‘normal_posts’" "GET ‘posts’" | Out-Null

# 2) POST 'normal_posts' -> POST 'posts'   (the plain "create a post" endpoint)
Replace-ExactRange "POST ‘normal_posts’" "POST ‘posts’" | Out-Null

# 3) POST 'normal_posts/comments' -> POST 'posts/comments'
Replace-ExactRange "POST ‘normal_posts/comments’" "POST ‘posts/comments’" | Out-Null

# 4) POST 'normal_posts/votes' -> POST 'posts/votes'
$r4 = Replace-ExactRange "POST ‘normal_posts/votes’" "POST ‘posts/votes’"

# The _GoBack bookmark (Word's "last edit location" marker) previously sat in
# its own empty paragraph after the "* id - the post id" params line; move it
# to mark this latest edit location instead, right after "POST '" and before
# "posts/votes'" (collapsed range at the start of the replaced text).
$prefixLen = "POST ‘".Length
$newEditPos = $d.Range($r4.Start + $prefixLen, $r4.Start + $prefixLen)
$d.Bookmarks.Add("_GoBack", $newEditPos) | Out-Null
